$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Propagate formatting (style, borders, alignment) from the last existing
# quarter column (AC) into the three new quarter columns before writing the
# new data, so e.g. the new header cells keep the bold/centered/bordered
# look of row 1, and the blank "section" rows stay blank cells.
$ws.Range("AC1:AC80").Copy($ws.Range("AD1:AD80"))
$ws.Range("AC1:AC80").Copy($ws.Range("AE1:AE80"))
$ws.Range("AC1:AC80").Copy($ws.Range("AF1:AF80"))

# Row 1: new period headers
$ws.Range("AD1").Value = "31/12/2023"
$ws.Range("AE1").Value = "31/03/2024"
$ws.Range("AF1").Value = "30/06/2024"

# Row 2
$ws.Range("AD2").Value = 14230176.768
$ws.Range("AE2").Value = 14746772.48
$ws.Range("AF2").Value = 16006472.704
# Row 3
$ws.Range("AD3").Value = 3346134.016
$ws.Range("AE3").Value = 3238871.04
$ws.Range("AF3").Value = 3706512.896
# Row 4
$ws.Range("AD4").Value = 1103432.96
$ws.Range("AE4").Value = 961820.032
$ws.Range("AF4").Value = 1434754.048
# Row 5
$ws.Range("AD5").Value = 1509880.064
$ws.Range("AE5").Value = 1375559.04
$ws.Range("AF5").Value = 1379828.992
# Row 6
$ws.Range("AD6").Value = 349007.008
$ws.Range("AE6").Value = 448632
$ws.Range("AF6").Value = 458787.008
# Row 7
$ws.Range("AD7").Value = 0
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 0
# Row 8
$ws.Range("AD8").Value = 0
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 0
# Row 9
$ws.Range("AD9").Value = 0
$ws.Range("AE9").Value = 0
$ws.Range("AF9").Value = 0
# Row 10
$ws.Range("AD10").Value = 0
$ws.Range("AE10").Value = 0
$ws.Range("AF10").Value = 0
# Row 11
$ws.Range("AD11").Value = 383814.016
$ws.Range("AE11").Value = 452860
$ws.Range("AF11").Value = 433143.008
# Row 12
$ws.Range("AD12").Value = 1130546.048
$ws.Range("AE12").Value = 1223560.96
$ws.Range("AF12").Value = 1284299.008
# Row 13
$ws.Range("AD13").Value = 33760
$ws.Range("AE13").Value = 33997
$ws.Range("AF13").Value = 35431
# Row 14
$ws.Range("AD14").Value = 0
$ws.Range("AE14").Value = 0
$ws.Range("AF14").Value = 0
# Row 15
$ws.Range("AD15").Value = 35452
$ws.Range("AE15").Value = 63828
$ws.Range("AF15").Value = 70333
# Row 16
$ws.Range("AD16").Value = 0
$ws.Range("AE16").Value = 0
$ws.Range("AF16").Value = 0
# Row 17
$ws.Range("AD17").Value = 0
$ws.Range("AE17").Value = 0
$ws.Range("AF17").Value = 0
# Row 18
$ws.Range("AD18").Value = 0
$ws.Range("AE18").Value = 0
$ws.Range("AF18").Value = 0
# Row 19
$ws.Range("AD19").Value = 798257.9840000001
$ws.Range("AE19").Value = 820297.024
$ws.Range("AF19").Value = 849803.008
# Row 20
$ws.Range("AD20").Value = 0
$ws.Range("AE20").Value = 0
$ws.Range("AF20").Value = 0
# Row 21
$ws.Range("AD21").Value = 0
$ws.Range("AE21").Value = 30384
$ws.Range("AF21").Value = 22581
# Row 22
$ws.Range("AD22").Value = 40940
$ws.Range("AE22").Value = 42612
$ws.Range("AF22").Value = 48384
# Row 23
$ws.Range("AD23").Value = 7799572.992
$ws.Range("AE23").Value = 8254796.8
$ws.Range("AF23").Value = 8892294.143999999
# Row 24
$ws.Range("AD24").Value = 1912984.064
$ws.Range("AE24").Value = 1986930.944
$ws.Range("AF24").Value = 2074983.04
# Row 25
$ws.Range("AD25").Value = 0
$ws.Range("AE25").Value = 0
$ws.Range("AF25").Value = 0
# Row 26
$ws.Range("AD26").Value = 14230176.768
$ws.Range("AE26").Value = 14746772.48
$ws.Range("AF26").Value = 16006472.704
# Row 27
$ws.Range("AD27").Value = 2244231.936
$ws.Range("AE27").Value = 2227118.08
$ws.Range("AF27").Value = 2264847.872
# Row 28
$ws.Range("AD28").Value = 96026
$ws.Range("AE28").Value = 108395
$ws.Range("AF28").Value = 124653
# Row 29
$ws.Range("AD29").Value = 399172
$ws.Range("AE29").Value = 341038.016
$ws.Range("AF29").Value = 379756
# Row 30
$ws.Range("AD30").Value = 216171.008
$ws.Range("AE30").Value = 238982
$ws.Range("AF30").Value = 222598
# Row 31
$ws.Range("AD31").Value = 594401.9840000001
$ws.Range("AE31").Value = 712769.9840000001
$ws.Range("AF31").Value = 638214.976
# Row 32
$ws.Range("AD32").Value = 2603
$ws.Range("AE32").Value = 1676
$ws.Range("AF32").Value = 1559
# Row 33
$ws.Range("AD33").Value = 0
$ws.Range("AE33").Value = 0
$ws.Range("AF33").Value = 47455
# Row 34
$ws.Range("AD34").Value = 935857.9840000001
$ws.Range("AE34").Value = 824257.024
$ws.Range("AF34").Value = 850612.008
# Row 35
$ws.Range("AD35").Value = 0
$ws.Range("AE35").Value = 0
$ws.Range("AF35").Value = 0
# Row 36
$ws.Range("AD36").Value = 0
$ws.Range("AE36").Value = 0
$ws.Range("AF36").Value = 0
# Row 37
$ws.Range("AD37").Value = 6890311.168
$ws.Range("AE37").Value = 7226664.96
$ws.Range("AF37").Value = 8386832.896
# Row 38
$ws.Range("AD38").Value = 3225497.088
$ws.Range("AE38").Value = 3361925.888
$ws.Range("AF38").Value = 4156131.072
# Row 39
$ws.Range("AD39").Value = 0
$ws.Range("AE39").Value = 0
$ws.Range("AF39").Value = 0
# Row 40
$ws.Range("AD40").Value = 3625760
$ws.Range("AE40").Value = 3824041.984
$ws.Range("AF40").Value = 4186278.912
# Row 41
$ws.Range("AD41").Value = 7719
$ws.Range("AE41").Value = 6766
$ws.Range("AF41").Value = 7153
# Row 42
$ws.Range("AD42").Value = 0
$ws.Range("AE42").Value = 0
$ws.Range("AF42").Value = 0
# Row 43
$ws.Range("AD43").Value = 31203
$ws.Range("AE43").Value = 33881
$ws.Range("AF43").Value = 32266
# Row 44
$ws.Range("AD44").Value = 0
$ws.Range("AE44").Value = 0
$ws.Range("AF44").Value = 0
# Row 45
$ws.Range("AD45").Value = 132
$ws.Range("AE45").Value = 50
$ws.Range("AF45").Value = 5004
# Row 46
$ws.Range("AD46").Value = 27372
$ws.Range("AE46").Value = 23228
$ws.Range("AF46").Value = 15175
# Row 47
$ws.Range("AD47").Value = 5068261.92
$ws.Range("AE47").Value = 5269760.928
$ws.Range("AF47").Value = 5339616.936
# Row 48
$ws.Range("AD48").Value = 2970443.008
$ws.Range("AE48").Value = 2970443.008
$ws.Range("AF48").Value = 2970443.008
# Row 49
$ws.Range("AD49").Value = 953420.992
$ws.Range("AE49").Value = 963302.0159999999
$ws.Range("AF49").Value = 836846.976
# Row 50
$ws.Range("AD50").Value = 0
$ws.Range("AE50").Value = 0
$ws.Range("AF50").Value = 0
# Row 51
$ws.Range("AD51").Value = 779011.968
$ws.Range("AE51").Value = 770553.9840000001
$ws.Range("AF51").Value = 770553.9840000001
# Row 52
$ws.Range("AD52").Value = -8458
$ws.Range("AE52").Value = 93225
$ws.Range("AF52").Value = 136291.008
# Row 53
$ws.Range("AD53").Value = 0
$ws.Range("AE53").Value = 0
$ws.Range("AF53").Value = 0
# Row 54
$ws.Range("AD54").Value = 0
$ws.Range("AE54").Value = 0
$ws.Range("AF54").Value = 0
# Row 55
$ws.Range("AD55").Value = 373844
$ws.Range("AE55").Value = 472236.992
$ws.Range("AF55").Value = 625481.9840000001
# Row 56
$ws.Range("AD56").Value = 0
$ws.Range("AE56").Value = 0
$ws.Range("AF56").Value = 0
# Row 59
$ws.Range("AD59").Value = 1131174.912
$ws.Range("AE59").Value = 1259888
$ws.Range("AF59").Value = 1357683.968
# Row 60
$ws.Range("AD60").Value = -675967.936
$ws.Range("AE60").Value = -736572.992
$ws.Range("AF60").Value = -802161.024
# Row 61
$ws.Range("AD61").Value = 455206.912
$ws.Range("AE61").Value = 523315.008
$ws.Range("AF61").Value = 555523.008
# Row 62
$ws.Range("AD62").Value = -97590.016
$ws.Range("AE62").Value = -102237
$ws.Range("AF62").Value = -107543
# Row 63
$ws.Range("AD63").Value = -134099.992
$ws.Range("AE63").Value = -132037
$ws.Range("AF63").Value = -136016.992
# Row 64
$ws.Range("AD64").Value = 0
$ws.Range("AE64").Value = 0
$ws.Range("AF64").Value = 0
# Row 65
$ws.Range("AD65").Value = 0
$ws.Range("AE65").Value = 0
$ws.Range("AF65").Value = 0
# Row 66
$ws.Range("AD66").Value = -2127.008
$ws.Range("AE66").Value = -10178
$ws.Range("AF66").Value = -2270
# Row 67
$ws.Range("AD67").Value = -2734
$ws.Range("AE67").Value = -1069
$ws.Range("AF67").Value = -547
# Row 68
$ws.Range("AD68").Value = -133513
$ws.Range("AE68").Value = -161239.008
$ws.Range("AF68").Value = -197738
# Row 69
$ws.Range("AD69").Value = 121207.992
$ws.Range("AE69").Value = 77829
$ws.Range("AF69").Value = 77158
# Row 70
$ws.Range("AD70").Value = -254720.992
$ws.Range("AE70").Value = -239068
$ws.Range("AF70").Value = -274896
# Row 74
$ws.Range("AD74").Value = 85142.992
$ws.Range("AE74").Value = 116555
$ws.Range("AF74").Value = 111408
# Row 75
$ws.Range("AD75").Value = -23640
$ws.Range("AE75").Value = -38123
$ws.Range("AF75").Value = -44473
# Row 76
$ws.Range("AD76").Value = 575097.9840000001
$ws.Range("AE76").Value = 16054
$ws.Range("AF76").Value = 27465
# Row 79
$ws.Range("AD79").Value = -3179
$ws.Range("AE79").Value = -1261
$ws.Range("AF79").Value = -1333
# Row 80
$ws.Range("AD80").Value = 633421.952
$ws.Range("AE80").Value = 93225
$ws.Range("AF80").Value = 93067
